# act tablas web jul25
# Updates the "120101" indicator workbook:
#  - Data sheet: adds a new "DINEM - MIDES" column (C) with the historic
#    series (1985-2018), relabels column B header to "MIDES-MEF-OPP" and
#    extends it with 2022/2023, reordering the year column (2023 -> 1985).
#  - Metadata sheet: rewrites "observaciones" text and adds a new
#    "actualizacion" / "Julio 2025" row.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("Data")
$md  = $wb.Worksheets.Item("Metadata")

# ---------------------------------------------------------------------
# Data sheet
# ---------------------------------------------------------------------

# Header row
$ws.Cells.Item(1,1).Value = "Fecha"
$ws.Cells.Item(1,2).Value = "MIDES-MEF-OPP"
$ws.Cells.Item(1,3).Value = "DINEM - MIDES"

# Year column (A) - new full order, 2023 down to 1985
$years = @(
    "2023","2022","2021","2020","2019","2018","2017","2016","2015","2014",
    "2013","2012","2011","2010","2009","2008","2007","2006","2005","2004",
    "2003","2002","2001","2000","1999","1998","1997","1996","1995","1994",
    "1993","1992","1991","1990","1989","1988","1987","1986","1985"
)

# Years are stored as text (not numbers) in the source workbook, same as
# the original "Fecha" column - force text storage so "2023" etc. isn't
# silently coerced to a numeric cell.
$ws.Range("A2:A40").NumberFormat = "@"

# Column B ("MIDES-MEF-OPP"): populated 2023..2005, blank 2004..1985
$colB = @(
    4.9,4.5,4.7,4.9,4.8,4.7,4.6,4.4,4.2,4.3,
    4.2,4.2,4.1,3.9,4.1,3.8,3.3,3.1,2.9,$null,
    $null,$null,$null,$null,$null,$null,$null,$null,$null,$null,
    $null,$null,$null,$null,$null,$null,$null,$null,$null
)

# Column C ("DINEM - MIDES"): blank 2023..2019, populated 2018..1985
$colC = @(
    $null,$null,$null,$null,$null,5.1,4.9,4.7,4.6,4.6,
    4.5,4.6,4.4,4.3,4.5,4.1,3.6,3.4,3.1,3.0,
    3.1,3.1,2.9,2.6,2.5,2.9,2.3,2.3,1.9,1.9,
    2.3,2.0,2.1,2.3,2.4,2.4,2.4,2.5,2.1
)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = 2 + $i

    $ws.Cells.Item($row, 1).Value = $years[$i]

    if ($colB[$i] -eq $null) {
        $ws.Cells.Item($row, 2).ClearContents()
    } else {
        $ws.Cells.Item($row, 2).Value = $colB[$i]
    }

    if ($colC[$i] -eq $null) {
        $ws.Cells.Item($row, 3).ClearContents()
    } else {
        $ws.Cells.Item($row, 3).Value = $colC[$i]
    }
}

# ---------------------------------------------------------------------
# Metadata sheet
# ---------------------------------------------------------------------

# The (previously empty) A1 placeholder cell is normalised to a single
# space, matching B1.
$md.Cells.Item(1, 1).Value = " "

$obs = "Las dos líneas representan metodologías ligeramente diferentes de cálculo. De acuerdo a lo informado en el Observatorio Social de MIDES, a partir del año 2016 se introdujo cambios en la metodología de estimación del Gasto Público Social producto de los cambios en la información brindada por el Presupuesto Nacional, lo cual llevó a trabajar en base al presupuesto por áreas programáticas (AP) de los incisos gubernamentales. El Gasto Público Social en Cultura y Deporte era considerado anteriormente bajo la denominación de Gasto Público Social No Convencional, definido como un subcomponente heterogéneo del GPS. La función Cultura y Deporte agrupa los gastos en museos, bibliotecas, organizaciones de prensa, servicios de televisión, deportes, y que antes también incluía otros conceptos que aludían a un aspecto multidisciplinario de los programas sociales. Se hizo una revisión de forma de dar consistencia en los conceptos para la serie desde 2015. La estimación siempre refiere a montos en pesos corrientes monto obligado intervenido por balance a partir de la información proporcionada mayoritariamente por Contaduría General de la Nación (CGN) del Ministerio de Economía y Finanzas (MEF). Para los años 2020 y 2021 se incluyen las erogaciones del fondo COVID destinadas a atender la emergencia sanitaria."

# Row 8 currently holds "observaciones" / "Sin observaciones" -> replace value
$md.Cells.Item(8, 2).Value = $obs

# Insert a new row 9 ("actualizacion" / "Julio 2025"), shifting the
# trailing "cita" rows down by one.
$md.Rows.Item(9).Insert()
$md.Cells.Item(9, 1).Value = "actualizacion"
$md.Cells.Item(9, 2).Value = "Julio 2025"
